# Updates the crypto price/volume/hour table to reflect the latest
# coinranking.com snapshot, and refreshes the top-20 coin ordering
# (rows 8-17) now that rankings have shifted.
#
# Commit: Updated symbol list on Mon Jan 23 21:00:09 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the text must be
# forced to stay text (Excel would otherwise auto-convert numeric-
# looking strings like "306.09", "1.52%" or "21" into real numbers).
$cellUpdates = @(
    @{ Cell = "D2"; Value = '306.09'; AsText = $true }
    @{ Cell = "E2"; Value = '1.52%'; AsText = $true }
    @{ Cell = "G2"; Value = '21'; AsText = $true }
    @{ Cell = "D3"; Value = '36.23'; AsText = $true }
    @{ Cell = "E3"; Value = '-1.65%'; AsText = $true }
    @{ Cell = "G3"; Value = '21'; AsText = $true }
    @{ Cell = "D4"; Value = '5.045'; AsText = $true }
    @{ Cell = "E4"; Value = '0.66%'; AsText = $true }
    @{ Cell = "G4"; Value = '21'; AsText = $true }
    @{ Cell = "D5"; Value = '0.07923'; AsText = $true }
    @{ Cell = "E5"; Value = '2.75%'; AsText = $true }
    @{ Cell = "G5"; Value = '21'; AsText = $true }
    @{ Cell = "D6"; Value = '2.289'; AsText = $true }
    @{ Cell = "E6"; Value = '11.80%'; AsText = $true }
    @{ Cell = "G6"; Value = '21'; AsText = $true }
    @{ Cell = "D7"; Value = '7.994'; AsText = $true }
    @{ Cell = "E7"; Value = '0.26%'; AsText = $true }
    @{ Cell = "G7"; Value = '21'; AsText = $true }
    @{ Cell = "B8"; Value = 'MXToken'; AsText = $false }
    @{ Cell = "C8"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; AsText = $false }
    @{ Cell = "D8"; Value = '0.9270'; AsText = $true }
    @{ Cell = "E8"; Value = '0.72%'; AsText = $true }
    @{ Cell = "G8"; Value = '21'; AsText = $true }
    @{ Cell = "B9"; Value = 'LiechtensteinCryptoassetsExchange'; AsText = $false }
    @{ Cell = "C9"; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; AsText = $false }
    @{ Cell = "D9"; Value = '0.09809'; AsText = $true }
    @{ Cell = "E9"; Value = '0.68%'; AsText = $true }
    @{ Cell = "G9"; Value = '21'; AsText = $true }
    @{ Cell = "B10"; Value = 'WazirX'; AsText = $false }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; AsText = $false }
    @{ Cell = "D10"; Value = '0.1861'; AsText = $true }
    @{ Cell = "E10"; Value = '0.41%'; AsText = $true }
    @{ Cell = "G10"; Value = '21'; AsText = $true }
    @{ Cell = "B11"; Value = 'MandalaExchangeToken'; AsText = $false }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; AsText = $false }
    @{ Cell = "D11"; Value = '0.08957'; AsText = $true }
    @{ Cell = "E11"; Value = '4.34%'; AsText = $true }
    @{ Cell = "G11"; Value = '21'; AsText = $true }
    @{ Cell = "B12"; Value = 'BitrueCoin'; AsText = $false }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; AsText = $false }
    @{ Cell = "D12"; Value = '0.03746'; AsText = $true }
    @{ Cell = "E12"; Value = '3.58%'; AsText = $true }
    @{ Cell = "G12"; Value = '21'; AsText = $true }
    @{ Cell = "B13"; Value = 'BitMartToken'; AsText = $false }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; AsText = $false }
    @{ Cell = "D13"; Value = '0.09914'; AsText = $true }
    @{ Cell = "E13"; Value = '-0.63%'; AsText = $true }
    @{ Cell = "G13"; Value = '21'; AsText = $true }
    @{ Cell = "B14"; Value = 'BitForexToken'; AsText = $false }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; AsText = $false }
    @{ Cell = "D14"; Value = '0.001433'; AsText = $true }
    @{ Cell = "E14"; Value = '-2.93%'; AsText = $true }
    @{ Cell = "G14"; Value = '21'; AsText = $true }
    @{ Cell = "B15"; Value = 'TigerCash'; AsText = $false }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; AsText = $false }
    @{ Cell = "D15"; Value = '0.005695'; AsText = $true }
    @{ Cell = "E15"; Value = '-0.99%'; AsText = $true }
    @{ Cell = "G15"; Value = '21'; AsText = $true }
    @{ Cell = "B16"; Value = 'LEO'; AsText = $false }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; AsText = $false }
    @{ Cell = "D16"; Value = '3.459'; AsText = $true }
    @{ Cell = "G16"; Value = '21'; AsText = $true }
    @{ Cell = "B17"; Value = 'GateToken'; AsText = $false }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; AsText = $false }
    @{ Cell = "D17"; Value = '4.150'; AsText = $true }
    @{ Cell = "E17"; Value = '2.48%'; AsText = $true }
    @{ Cell = "G17"; Value = '21'; AsText = $true }
    @{ Cell = "E18"; Value = '4.35%'; AsText = $true }
    @{ Cell = "G18"; Value = '21'; AsText = $true }
    @{ Cell = "D19"; Value = '0.3368'; AsText = $true }
    @{ Cell = "E19"; Value = '-1.79%'; AsText = $true }
    @{ Cell = "G19"; Value = '21'; AsText = $true }
    @{ Cell = "D20"; Value = '0.1320'; AsText = $true }
    @{ Cell = "E20"; Value = '-1.06%'; AsText = $true }
    @{ Cell = "G20"; Value = '21'; AsText = $true }
    @{ Cell = "D21"; Value = '5.069'; AsText = $true }
    @{ Cell = "E21"; Value = '1.89%'; AsText = $true }
    @{ Cell = "G21"; Value = '21'; AsText = $true }
    @{ Cell = "D22"; Value = '0.2248'; AsText = $true }
    @{ Cell = "E22"; Value = '1.50%'; AsText = $true }
    @{ Cell = "G22"; Value = '21'; AsText = $true }
    @{ Cell = "D23"; Value = '0.04577'; AsText = $true }
    @{ Cell = "E23"; Value = '-0.77%'; AsText = $true }
    @{ Cell = "G23"; Value = '21'; AsText = $true }
    @{ Cell = "D24"; Value = '0.001233'; AsText = $true }
    @{ Cell = "E24"; Value = '-0.71%'; AsText = $true }
    @{ Cell = "G24"; Value = '21'; AsText = $true }
    @{ Cell = "D25"; Value = '0.004778'; AsText = $true }
    @{ Cell = "E25"; Value = '-6.29%'; AsText = $true }
    @{ Cell = "G25"; Value = '21'; AsText = $true }
    @{ Cell = "D26"; Value = '0.0001300'; AsText = $true }
    @{ Cell = "E26"; Value = '-8.04%'; AsText = $true }
    @{ Cell = "G26"; Value = '21'; AsText = $true }
    @{ Cell = "G27"; Value = '21'; AsText = $true }
    @{ Cell = "G28"; Value = '21'; AsText = $true }
    @{ Cell = "G29"; Value = '21'; AsText = $true }
    @{ Cell = "G30"; Value = '21'; AsText = $true }
    @{ Cell = "G31"; Value = '21'; AsText = $true }
    @{ Cell = "G32"; Value = '21'; AsText = $true }
    @{ Cell = "G33"; Value = '21'; AsText = $true }
    @{ Cell = "G34"; Value = '21'; AsText = $true }
    @{ Cell = "G35"; Value = '21'; AsText = $true }
    @{ Cell = "G36"; Value = '21'; AsText = $true }
    @{ Cell = "G37"; Value = '21'; AsText = $true }
    @{ Cell = "G38"; Value = '21'; AsText = $true }
    @{ Cell = "D39"; Value = '0.01921'; AsText = $true }
    @{ Cell = "E39"; Value = '10.30%'; AsText = $true }
    @{ Cell = "G39"; Value = '21'; AsText = $true }
    @{ Cell = "D40"; Value = '0.04908'; AsText = $true }
    @{ Cell = "E40"; Value = '6.32%'; AsText = $true }
    @{ Cell = "G40"; Value = '21'; AsText = $true }
    @{ Cell = "D41"; Value = '0.007784'; AsText = $true }
    @{ Cell = "E41"; Value = '1.13%'; AsText = $true }
    @{ Cell = "G41"; Value = '21'; AsText = $true }
    @{ Cell = "D42"; Value = '0.1392'; AsText = $true }
    @{ Cell = "E42"; Value = '-0.01%'; AsText = $true }
    @{ Cell = "G42"; Value = '21'; AsText = $true }
    @{ Cell = "D43"; Value = '0.007802'; AsText = $true }
    @{ Cell = "E43"; Value = '-2.47%'; AsText = $true }
    @{ Cell = "G43"; Value = '21'; AsText = $true }
    @{ Cell = "D44"; Value = '0.002118'; AsText = $true }
    @{ Cell = "E44"; Value = '-2.61%'; AsText = $true }
    @{ Cell = "G44"; Value = '21'; AsText = $true }
    @{ Cell = "E45"; Value = '15.45%'; AsText = $true }
    @{ Cell = "G45"; Value = '21'; AsText = $true }
    @{ Cell = "D46"; Value = '0.00006141'; AsText = $true }
    @{ Cell = "E46"; Value = '-2.79%'; AsText = $true }
    @{ Cell = "G46"; Value = '21'; AsText = $true }
    @{ Cell = "D47"; Value = '0.00000000750'; AsText = $true }
    @{ Cell = "E47"; Value = '-0.94%'; AsText = $true }
    @{ Cell = "G47"; Value = '21'; AsText = $true }
    @{ Cell = "D48"; Value = '51.77'; AsText = $true }
    @{ Cell = "E48"; Value = '52.96%'; AsText = $true }
    @{ Cell = "G48"; Value = '21'; AsText = $true }
    @{ Cell = "D49"; Value = '0.001800'; AsText = $true }
    @{ Cell = "E49"; Value = '-10.83%'; AsText = $true }
    @{ Cell = "G49"; Value = '21'; AsText = $true }
    @{ Cell = "D50"; Value = '0.00002100'; AsText = $true }
    @{ Cell = "E50"; Value = '-0.94%'; AsText = $true }
    @{ Cell = "G50"; Value = '21'; AsText = $true }
    @{ Cell = "D51"; Value = '0.0002000'; AsText = $true }
    @{ Cell = "E51"; Value = '-0.94%'; AsText = $true }
    @{ Cell = "G51"; Value = '21'; AsText = $true }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    if ($update.AsText) {
        # Keep the original text representation (e.g. trailing zeros,
        # percent signs) instead of Excel's automatic number parsing.
        $range.NumberFormat = "@"
    }
    $range.Value = $update.Value
}
